$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as plain text in the source sheet (e.g.
# thousands-dot-separated strings like "66.844.94"). Setting .Value directly would
# let Excel auto-coerce plain-decimal-looking strings (e.g. "585.85") into numbers,
# so force the cell to Text format first, write the literal string, then restore
# General so the cell keeps behaving like any other unformatted cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.844.94"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -1.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.608.18"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -1.75%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.85"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -1.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.19"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -1.53%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.527"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -3.72%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.606.44"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -1.81%  "

$ws.Range("E10").Value = "  -3.63%  "

$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.364"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.72%  "

$ws.Range("E13").Value = "  -1.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.11"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -3.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.086.73"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -1.64%  "

$ws.Range("E16").Value = "  -3.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.843.05"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -1.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.615.83"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -1.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.63"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -4.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.76"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -5.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "353.92"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -2.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.25"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -3.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.60"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -4.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.46"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -5.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.89"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -6.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "69.18"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -2.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.745.69"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -1.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0986"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -4.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "538.02"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -3.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.08"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +0.20%  "

$ws.Range("E33").Value = "  -4.81%  "

$ws.Range("E34").Value = "  -3.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.132"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -2.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("E37").Value = "  -5.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.95"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.83"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -2.90%  "

$ws.Range("E40").Value = "  -2.89%  "

$ws.Range("E41").Value = "  +1.83%  "

$ws.Range("E42").Value = "  -1.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.10"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -4.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.39"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -6.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0291"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -2.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "150.24"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -2.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.572"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -4.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.73"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -3.81%  "

$ws.Range("E50").Value = "  -2.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0766"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -1.71%  "
